$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "[Publication]"
$ws.Range("A27").Font.Bold = $true
$ws.Range("A28").Value = "#DOI"
$ws.Range("A29").Value = "abc/def"

$ws.Range("A28").Select()
